# Reset vouchers.xlsx: clear all voucher data, keep only a fresh header row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the existing data rows (2-4), leaving just the header row.
$ws.Rows("2:4").Delete()

# Rewrite the header row with the new column set/order.
$ws.Range("A1").Value = "voucher_code"
$ws.Range("B1").Value = "phone"
$ws.Range("C1").Value = "value"
$ws.Range("D1").Value = "issued_ts"
$ws.Range("E1").Value = "redeemed_ts"

# Drop the two trailing columns (old "redeemed" + "redeemed_ts") entirely
# so the used range shrinks back down to A1:E1.
$ws.Range("F1:G1").EntireColumn.Delete()
